$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update DATA_EXTRACCIO (column H) timestamps for all data rows ---
$ws.Range("H2:H226").NumberFormat = "@"

$hValues = @{
    2 = "2026-02-20 08:15:27"
    3 = "2026-02-20 08:15:28"
    4 = "2026-02-20 08:15:28"
    5 = "2026-02-20 08:15:28"
    6 = "2026-02-20 08:15:28"
    7 = "2026-02-20 08:15:31"
    8 = "2026-02-20 08:15:32"
    9 = "2026-02-20 08:15:32"
    10 = "2026-02-20 08:15:32"
    11 = "2026-02-20 08:15:32"
    12 = "2026-02-20 08:15:35"
    13 = "2026-02-20 08:15:36"
    14 = "2026-02-20 08:15:36"
    15 = "2026-02-20 08:15:36"
    16 = "2026-02-20 08:15:36"
    17 = "2026-02-20 08:15:39"
    18 = "2026-02-20 08:15:40"
    19 = "2026-02-20 08:15:40"
    20 = "2026-02-20 08:15:40"
    21 = "2026-02-20 08:15:40"
    22 = "2026-02-20 08:15:43"
    23 = "2026-02-20 08:15:44"
    24 = "2026-02-20 08:15:44"
    25 = "2026-02-20 08:15:44"
    26 = "2026-02-20 08:15:44"
    27 = "2026-02-20 08:15:47"
    28 = "2026-02-20 08:15:48"
    29 = "2026-02-20 08:15:48"
    30 = "2026-02-20 08:15:48"
    31 = "2026-02-20 08:15:48"
    32 = "2026-02-20 08:15:51"
    33 = "2026-02-20 08:15:52"
    34 = "2026-02-20 08:15:52"
    35 = "2026-02-20 08:15:52"
    36 = "2026-02-20 08:15:52"
    37 = "2026-02-20 08:15:55"
    38 = "2026-02-20 08:15:56"
    39 = "2026-02-20 08:15:56"
    40 = "2026-02-20 08:15:56"
    41 = "2026-02-20 08:15:56"
    42 = "2026-02-20 08:15:58"
    43 = "2026-02-20 08:16:00"
    44 = "2026-02-20 08:16:00"
    45 = "2026-02-20 08:16:00"
    46 = "2026-02-20 08:16:00"
    47 = "2026-02-20 08:16:02"
    48 = "2026-02-20 08:16:03"
    49 = "2026-02-20 08:16:03"
    50 = "2026-02-20 08:16:03"
    51 = "2026-02-20 08:16:03"
    52 = "2026-02-20 08:16:06"
    53 = "2026-02-20 08:16:07"
    54 = "2026-02-20 08:16:07"
    55 = "2026-02-20 08:16:07"
    56 = "2026-02-20 08:16:07"
    57 = "2026-02-20 08:16:10"
    58 = "2026-02-20 08:16:11"
    59 = "2026-02-20 08:16:11"
    60 = "2026-02-20 08:16:11"
    61 = "2026-02-20 08:16:11"
    62 = "2026-02-20 08:16:14"
    63 = "2026-02-20 08:16:15"
    64 = "2026-02-20 08:16:15"
    65 = "2026-02-20 08:16:15"
    66 = "2026-02-20 08:16:15"
    67 = "2026-02-20 08:16:17"
    68 = "2026-02-20 08:16:19"
    69 = "2026-02-20 08:16:19"
    70 = "2026-02-20 08:16:19"
    71 = "2026-02-20 08:16:19"
    72 = "2026-02-20 08:16:21"
    73 = "2026-02-20 08:16:23"
    74 = "2026-02-20 08:16:23"
    75 = "2026-02-20 08:16:23"
    76 = "2026-02-20 08:16:23"
    77 = "2026-02-20 08:16:25"
    78 = "2026-02-20 08:16:26"
    79 = "2026-02-20 08:16:26"
    80 = "2026-02-20 08:16:26"
    81 = "2026-02-20 08:16:26"
    82 = "2026-02-20 08:16:29"
    83 = "2026-02-20 08:16:30"
    84 = "2026-02-20 08:16:30"
    85 = "2026-02-20 08:16:30"
    86 = "2026-02-20 08:16:30"
    87 = "2026-02-20 08:16:32"
    88 = "2026-02-20 08:16:34"
    89 = "2026-02-20 08:16:34"
    90 = "2026-02-20 08:16:34"
    91 = "2026-02-20 08:16:34"
    92 = "2026-02-20 08:16:36"
    93 = "2026-02-20 08:16:38"
    94 = "2026-02-20 08:16:38"
    95 = "2026-02-20 08:16:38"
    96 = "2026-02-20 08:16:38"
    97 = "2026-02-20 08:16:40"
    98 = "2026-02-20 08:16:41"
    99 = "2026-02-20 08:16:41"
    100 = "2026-02-20 08:16:41"
    101 = "2026-02-20 08:16:41"
    102 = "2026-02-20 08:16:44"
    103 = "2026-02-20 08:16:45"
    104 = "2026-02-20 08:16:45"
    105 = "2026-02-20 08:16:45"
    106 = "2026-02-20 08:16:45"
    107 = "2026-02-20 08:16:48"
    108 = "2026-02-20 08:16:49"
    109 = "2026-02-20 08:16:49"
    110 = "2026-02-20 08:16:49"
    111 = "2026-02-20 08:16:49"
    112 = "2026-02-20 08:16:51"
    113 = "2026-02-20 08:16:53"
    114 = "2026-02-20 08:16:53"
    115 = "2026-02-20 08:16:53"
    116 = "2026-02-20 08:16:53"
    117 = "2026-02-20 08:16:55"
    118 = "2026-02-20 08:16:57"
    119 = "2026-02-20 08:16:57"
    120 = "2026-02-20 08:16:57"
    121 = "2026-02-20 08:16:57"
    122 = "2026-02-20 08:16:59"
    123 = "2026-02-20 08:17:01"
    124 = "2026-02-20 08:17:01"
    125 = "2026-02-20 08:17:01"
    126 = "2026-02-20 08:17:01"
    127 = "2026-02-20 08:17:03"
    128 = "2026-02-20 08:17:04"
    129 = "2026-02-20 08:17:04"
    130 = "2026-02-20 08:17:04"
    131 = "2026-02-20 08:17:04"
    132 = "2026-02-20 08:17:06"
    133 = "2026-02-20 08:17:08"
    134 = "2026-02-20 08:17:08"
    135 = "2026-02-20 08:17:08"
    136 = "2026-02-20 08:17:08"
    137 = "2026-02-20 08:17:10"
    138 = "2026-02-20 08:17:12"
    139 = "2026-02-20 08:17:12"
    140 = "2026-02-20 08:17:12"
    141 = "2026-02-20 08:17:12"
    142 = "2026-02-20 08:17:14"
    143 = "2026-02-20 08:17:15"
    144 = "2026-02-20 08:17:15"
    145 = "2026-02-20 08:17:15"
    146 = "2026-02-20 08:17:15"
    147 = "2026-02-20 08:17:18"
    148 = "2026-02-20 08:17:19"
    149 = "2026-02-20 08:17:19"
    150 = "2026-02-20 08:17:19"
    151 = "2026-02-20 08:17:19"
    152 = "2026-02-20 08:17:22"
    153 = "2026-02-20 08:17:23"
    154 = "2026-02-20 08:17:23"
    155 = "2026-02-20 08:17:23"
    156 = "2026-02-20 08:17:23"
    157 = "2026-02-20 08:17:25"
    158 = "2026-02-20 08:17:27"
    159 = "2026-02-20 08:17:27"
    160 = "2026-02-20 08:17:27"
    161 = "2026-02-20 08:17:27"
    162 = "2026-02-20 08:17:29"
    163 = "2026-02-20 08:17:31"
    164 = "2026-02-20 08:17:31"
    165 = "2026-02-20 08:17:31"
    166 = "2026-02-20 08:17:31"
    167 = "2026-02-20 08:17:33"
    168 = "2026-02-20 08:17:35"
    169 = "2026-02-20 08:17:35"
    170 = "2026-02-20 08:17:35"
    171 = "2026-02-20 08:17:35"
    172 = "2026-02-20 08:17:37"
    173 = "2026-02-20 08:17:39"
    174 = "2026-02-20 08:17:39"
    175 = "2026-02-20 08:17:39"
    176 = "2026-02-20 08:17:39"
    177 = "2026-02-20 08:17:41"
    178 = "2026-02-20 08:17:42"
    179 = "2026-02-20 08:17:42"
    180 = "2026-02-20 08:17:42"
    181 = "2026-02-20 08:17:42"
    182 = "2026-02-20 08:17:45"
    183 = "2026-02-20 08:17:46"
    184 = "2026-02-20 08:17:46"
    185 = "2026-02-20 08:17:46"
    186 = "2026-02-20 08:17:46"
    187 = "2026-02-20 08:17:49"
    188 = "2026-02-20 08:17:50"
    189 = "2026-02-20 08:17:50"
    190 = "2026-02-20 08:17:50"
    191 = "2026-02-20 08:17:50"
    192 = "2026-02-20 08:17:52"
    193 = "2026-02-20 08:17:54"
    194 = "2026-02-20 08:17:54"
    195 = "2026-02-20 08:17:54"
    196 = "2026-02-20 08:17:54"
    197 = "2026-02-20 08:17:56"
    198 = "2026-02-20 08:17:58"
    199 = "2026-02-20 08:17:58"
    200 = "2026-02-20 08:17:58"
    201 = "2026-02-20 08:17:58"
    202 = "2026-02-20 08:18:00"
    203 = "2026-02-20 08:18:01"
    204 = "2026-02-20 08:18:01"
    205 = "2026-02-20 08:18:01"
    206 = "2026-02-20 08:18:01"
    207 = "2026-02-20 08:18:04"
    208 = "2026-02-20 08:18:05"
    209 = "2026-02-20 08:18:05"
    210 = "2026-02-20 08:18:05"
    211 = "2026-02-20 08:18:05"
    212 = "2026-02-20 08:18:08"
    213 = "2026-02-20 08:18:09"
    214 = "2026-02-20 08:18:09"
    215 = "2026-02-20 08:18:09"
    216 = "2026-02-20 08:18:09"
    217 = "2026-02-20 08:18:12"
    218 = "2026-02-20 08:18:13"
    219 = "2026-02-20 08:18:13"
    220 = "2026-02-20 08:18:13"
    221 = "2026-02-20 08:18:13"
    222 = "2026-02-20 08:18:16"
    223 = "2026-02-20 08:18:17"
    224 = "2026-02-20 08:18:17"
    225 = "2026-02-20 08:18:17"
    226 = "2026-02-20 08:18:17"
}

foreach ($r in $hValues.Keys) {
    $ws.Cells.Item($r, 8).Value = $hValues[$r]
}

# --- Update the 7 "current period" rows (07:00-07:30 -> 07:30-08:00) with new data ---

# Row 2
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "07:30 - 08:00"
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "116"
$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "66"
$ws.Range("T2").NumberFormat = "@"
$ws.Range("T2").Value = "07:30 - 08:00"
$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value = "53"
$ws.Range("V2").NumberFormat = "@"
$ws.Range("V2").Value = "-0.8"
$ws.Range("X2").NumberFormat = "@"
$ws.Range("X2").Value = "-0.6"
$ws.Range("AE2").NumberFormat = "@"
$ws.Range("AE2").Value = "07:30 - 08:00"
$ws.Range("AF2").NumberFormat = "@"
$ws.Range("AF2").Value = "-0.8"
$ws.Range("AG2").NumberFormat = "@"
$ws.Range("AG2").Value = "-0.6"
$ws.Range("AI2").NumberFormat = "@"
$ws.Range("AI2").Value = "66"
$ws.Range("AK2").NumberFormat = "@"
$ws.Range("AK2").Value = "116"
$ws.Range("AL2").NumberFormat = "@"
$ws.Range("AL2").Value = "53"
$ws.Range("AM2").NumberFormat = "@"
$ws.Range("AM2").Value = "07:30 - 08:00"
$ws.Range("AN2").NumberFormat = "@"
$ws.Range("AN2").Value = "-0.8"
$ws.Range("AO2").NumberFormat = "@"
$ws.Range("AO2").Value = "-0.6"
$ws.Range("AQ2").NumberFormat = "@"
$ws.Range("AQ2").Value = "66"
$ws.Range("AS2").NumberFormat = "@"
$ws.Range("AS2").Value = "116"
$ws.Range("AT2").NumberFormat = "@"
$ws.Range("AT2").Value = "53"

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "07:30 - 08:00"
$ws.Range("O7").NumberFormat = "@"
$ws.Range("O7").Value = "254"
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "262"
$ws.Range("T7").NumberFormat = "@"
$ws.Range("T7").Value = "07:30 - 08:00"
$ws.Range("U7").NumberFormat = "@"
$ws.Range("U7").Value = "82"
$ws.Range("V7").NumberFormat = "@"
$ws.Range("V7").Value = "-6.1"
$ws.Range("W7").NumberFormat = "@"
$ws.Range("W7").Value = "-6.2"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "25.2"
$ws.Range("AD7").NumberFormat = "@"
$ws.Range("AD7").Value = "35.6"
$ws.Range("AE7").NumberFormat = "@"
$ws.Range("AE7").Value = "07:30 - 08:00"
$ws.Range("AF7").NumberFormat = "@"
$ws.Range("AF7").Value = "-6.1"
$ws.Range("AH7").NumberFormat = "@"
$ws.Range("AH7").Value = "-6.2"
$ws.Range("AK7").NumberFormat = "@"
$ws.Range("AK7").Value = "262"
$ws.Range("AL7").NumberFormat = "@"
$ws.Range("AL7").Value = "25.2"
$ws.Range("AM7").NumberFormat = "@"
$ws.Range("AM7").Value = "07:30 - 08:00"
$ws.Range("AN7").NumberFormat = "@"
$ws.Range("AN7").Value = "-6.1"
$ws.Range("AP7").NumberFormat = "@"
$ws.Range("AP7").Value = "-6.2"
$ws.Range("AS7").NumberFormat = "@"
$ws.Range("AS7").Value = "262"
$ws.Range("AT7").NumberFormat = "@"
$ws.Range("AT7").Value = "25.2"
$ws.Range("AU7").NumberFormat = "@"
$ws.Range("AU7").Value = "254"
$ws.Range("AV7").NumberFormat = "@"
$ws.Range("AV7").Value = "35.6"
$ws.Range("AW7").NumberFormat = "@"
$ws.Range("AW7").Value = "82"
$ws.Range("AX7").NumberFormat = "@"
$ws.Range("AX7").Value = "254"
$ws.Range("AY7").NumberFormat = "@"
$ws.Range("AY7").Value = "35.6"
$ws.Range("AZ7").NumberFormat = "@"
$ws.Range("AZ7").Value = "82"

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "07:30 - 08:00"
$ws.Range("M12").NumberFormat = "@"
$ws.Range("M12").Value = "341"
$ws.Range("Q12").NumberFormat = "@"
$ws.Range("Q12").Value = "56"
$ws.Range("R12").NumberFormat = "@"
$ws.Range("R12").Value = "1023.2"
$ws.Range("T12").NumberFormat = "@"
$ws.Range("T12").Value = "07:30 - 08:00"
$ws.Range("U12").NumberFormat = "@"
$ws.Range("U12").Value = "198"
$ws.Range("V12").NumberFormat = "@"
$ws.Range("V12").Value = "8.2"
$ws.Range("W12").NumberFormat = "@"
$ws.Range("W12").Value = "7.9"
$ws.Range("X12").NumberFormat = "@"
$ws.Range("X12").Value = "8.7"
$ws.Range("Y12").NumberFormat = "@"
$ws.Range("Y12").Value = "11.2"
$ws.Range("AB12").NumberFormat = "@"
$ws.Range("AB12").Value = "25.9"
$ws.Range("AE12").NumberFormat = "@"
$ws.Range("AE12").Value = "07:30 - 08:00"
$ws.Range("AF12").NumberFormat = "@"
$ws.Range("AF12").Value = "8.2"
$ws.Range("AG12").NumberFormat = "@"
$ws.Range("AG12").Value = "8.7"
$ws.Range("AH12").NumberFormat = "@"
$ws.Range("AH12").Value = "7.9"
$ws.Range("AI12").NumberFormat = "@"
$ws.Range("AI12").Value = "56"
$ws.Range("AK12").NumberFormat = "@"
$ws.Range("AK12").Value = "11.2"
$ws.Range("AL12").NumberFormat = "@"
$ws.Range("AL12").Value = "341"
$ws.Range("AM12").NumberFormat = "@"
$ws.Range("AM12").Value = "07:30 - 08:00"
$ws.Range("AN12").NumberFormat = "@"
$ws.Range("AN12").Value = "8.2"
$ws.Range("AO12").NumberFormat = "@"
$ws.Range("AO12").Value = "8.7"
$ws.Range("AP12").NumberFormat = "@"
$ws.Range("AP12").Value = "7.9"
$ws.Range("AQ12").NumberFormat = "@"
$ws.Range("AQ12").Value = "56"
$ws.Range("AS12").NumberFormat = "@"
$ws.Range("AS12").Value = "11.2"
$ws.Range("AT12").NumberFormat = "@"
$ws.Range("AT12").Value = "341"
$ws.Range("AU12").NumberFormat = "@"
$ws.Range("AU12").Value = "25.9"
$ws.Range("AV12").NumberFormat = "@"
$ws.Range("AV12").Value = "1023.2"
$ws.Range("AW12").NumberFormat = "@"
$ws.Range("AW12").Value = "198"
$ws.Range("AX12").NumberFormat = "@"
$ws.Range("AX12").Value = "25.9"
$ws.Range("AY12").NumberFormat = "@"
$ws.Range("AY12").Value = "1023.2"
$ws.Range("AZ12").NumberFormat = "@"
$ws.Range("AZ12").Value = "198"

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "07:30 - 08:00"
$ws.Range("O17").NumberFormat = "@"
$ws.Range("O17").Value = "316"
$ws.Range("P17").NumberFormat = "@"
$ws.Range("P17").Value = "130"
$ws.Range("Q17").NumberFormat = "@"
$ws.Range("Q17").Value = "92"
$ws.Range("T17").NumberFormat = "@"
$ws.Range("T17").Value = "07:30 - 08:00"
$ws.Range("U17").NumberFormat = "@"
$ws.Range("U17").Value = "53"
$ws.Range("V17").NumberFormat = "@"
$ws.Range("V17").Value = "-5.7"
$ws.Range("W17").NumberFormat = "@"
$ws.Range("W17").Value = "-5.8"
$ws.Range("X17").NumberFormat = "@"
$ws.Range("X17").Value = "-5.6"
$ws.Range("AD17").NumberFormat = "@"
$ws.Range("AD17").Value = "15.8"
$ws.Range("AE17").NumberFormat = "@"
$ws.Range("AE17").Value = "07:30 - 08:00"
$ws.Range("AF17").NumberFormat = "@"
$ws.Range("AF17").Value = "-5.7"
$ws.Range("AG17").NumberFormat = "@"
$ws.Range("AG17").Value = "-5.6"
$ws.Range("AH17").NumberFormat = "@"
$ws.Range("AH17").Value = "-5.8"
$ws.Range("AI17").NumberFormat = "@"
$ws.Range("AI17").Value = "92"
$ws.Range("AK17").NumberFormat = "@"
$ws.Range("AK17").Value = "130"
$ws.Range("AM17").NumberFormat = "@"
$ws.Range("AM17").Value = "07:30 - 08:00"
$ws.Range("AN17").NumberFormat = "@"
$ws.Range("AN17").Value = "-5.7"
$ws.Range("AO17").NumberFormat = "@"
$ws.Range("AO17").Value = "-5.6"
$ws.Range("AP17").NumberFormat = "@"
$ws.Range("AP17").Value = "-5.8"
$ws.Range("AQ17").NumberFormat = "@"
$ws.Range("AQ17").Value = "92"
$ws.Range("AS17").NumberFormat = "@"
$ws.Range("AS17").Value = "130"
$ws.Range("AU17").NumberFormat = "@"
$ws.Range("AU17").Value = "316"
$ws.Range("AV17").NumberFormat = "@"
$ws.Range("AV17").Value = "15.8"
$ws.Range("AW17").NumberFormat = "@"
$ws.Range("AW17").Value = "53"
$ws.Range("AX17").NumberFormat = "@"
$ws.Range("AX17").Value = "316"
$ws.Range("AY17").NumberFormat = "@"
$ws.Range("AY17").Value = "15.8"
$ws.Range("AZ17").NumberFormat = "@"
$ws.Range("AZ17").Value = "53"

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "07:30 - 08:00"
$ws.Range("M22").NumberFormat = "@"
$ws.Range("M22").Value = "319"
$ws.Range("Q22").NumberFormat = "@"
$ws.Range("Q22").Value = "68"
$ws.Range("R22").NumberFormat = "@"
$ws.Range("R22").Value = "1023.2"
$ws.Range("T22").NumberFormat = "@"
$ws.Range("T22").Value = "07:30 - 08:00"
$ws.Range("U22").NumberFormat = "@"
$ws.Range("U22").Value = "152"
$ws.Range("V22").NumberFormat = "@"
$ws.Range("V22").Value = "7.4"
$ws.Range("W22").NumberFormat = "@"
$ws.Range("W22").Value = "6.1"
$ws.Range("X22").NumberFormat = "@"
$ws.Range("X22").Value = "8.8"
$ws.Range("Y22").NumberFormat = "@"
$ws.Range("Y22").Value = "5.8"
$ws.Range("AB22").NumberFormat = "@"
$ws.Range("AB22").Value = "13.0"
$ws.Range("AE22").NumberFormat = "@"
$ws.Range("AE22").Value = "07:30 - 08:00"
$ws.Range("AF22").NumberFormat = "@"
$ws.Range("AF22").Value = "7.4"
$ws.Range("AG22").NumberFormat = "@"
$ws.Range("AG22").Value = "8.8"
$ws.Range("AH22").NumberFormat = "@"
$ws.Range("AH22").Value = "6.1"
$ws.Range("AI22").NumberFormat = "@"
$ws.Range("AI22").Value = "68"
$ws.Range("AK22").NumberFormat = "@"
$ws.Range("AK22").Value = "5.8"
$ws.Range("AL22").NumberFormat = "@"
$ws.Range("AL22").Value = "319"
$ws.Range("AM22").NumberFormat = "@"
$ws.Range("AM22").Value = "07:30 - 08:00"
$ws.Range("AN22").NumberFormat = "@"
$ws.Range("AN22").Value = "7.4"
$ws.Range("AO22").NumberFormat = "@"
$ws.Range("AO22").Value = "8.8"
$ws.Range("AP22").NumberFormat = "@"
$ws.Range("AP22").Value = "6.1"
$ws.Range("AQ22").NumberFormat = "@"
$ws.Range("AQ22").Value = "68"
$ws.Range("AS22").NumberFormat = "@"
$ws.Range("AS22").Value = "5.8"
$ws.Range("AT22").NumberFormat = "@"
$ws.Range("AT22").Value = "319"
$ws.Range("AU22").NumberFormat = "@"
$ws.Range("AU22").Value = "13.0"
$ws.Range("AV22").NumberFormat = "@"
$ws.Range("AV22").Value = "1023.2"
$ws.Range("AW22").NumberFormat = "@"
$ws.Range("AW22").Value = "152"
$ws.Range("AX22").NumberFormat = "@"
$ws.Range("AX22").Value = "13.0"
$ws.Range("AY22").NumberFormat = "@"
$ws.Range("AY22").Value = "1023.2"
$ws.Range("AZ22").NumberFormat = "@"
$ws.Range("AZ22").Value = "152"

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "07:30 - 08:00"
$ws.Range("M27").NumberFormat = "@"
$ws.Range("M27").Value = "336"
$ws.Range("Q27").NumberFormat = "@"
$ws.Range("Q27").Value = "45"
$ws.Range("R27").NumberFormat = "@"
$ws.Range("R27").Value = "1022.6"
$ws.Range("T27").NumberFormat = "@"
$ws.Range("T27").Value = "07:30 - 08:00"
$ws.Range("U27").NumberFormat = "@"
$ws.Range("U27").Value = "132"
$ws.Range("V27").NumberFormat = "@"
$ws.Range("V27").Value = "11.7"
$ws.Range("X27").NumberFormat = "@"
$ws.Range("X27").Value = "12.1"
$ws.Range("Y27").NumberFormat = "@"
$ws.Range("Y27").Value = "5.0"
$ws.Range("AB27").NumberFormat = "@"
$ws.Range("AB27").Value = "17.3"
$ws.Range("AE27").NumberFormat = "@"
$ws.Range("AE27").Value = "07:30 - 08:00"
$ws.Range("AF27").NumberFormat = "@"
$ws.Range("AF27").Value = "11.7"
$ws.Range("AG27").NumberFormat = "@"
$ws.Range("AG27").Value = "12.1"
$ws.Range("AI27").NumberFormat = "@"
$ws.Range("AI27").Value = "45"
$ws.Range("AK27").NumberFormat = "@"
$ws.Range("AK27").Value = "5.0"
$ws.Range("AL27").NumberFormat = "@"
$ws.Range("AL27").Value = "336"
$ws.Range("AM27").NumberFormat = "@"
$ws.Range("AM27").Value = "07:30 - 08:00"
$ws.Range("AN27").NumberFormat = "@"
$ws.Range("AN27").Value = "11.7"
$ws.Range("AO27").NumberFormat = "@"
$ws.Range("AO27").Value = "12.1"
$ws.Range("AQ27").NumberFormat = "@"
$ws.Range("AQ27").Value = "45"
$ws.Range("AS27").NumberFormat = "@"
$ws.Range("AS27").Value = "5.0"
$ws.Range("AT27").NumberFormat = "@"
$ws.Range("AT27").Value = "336"
$ws.Range("AU27").NumberFormat = "@"
$ws.Range("AU27").Value = "17.3"
$ws.Range("AV27").NumberFormat = "@"
$ws.Range("AV27").Value = "1022.6"
$ws.Range("AW27").NumberFormat = "@"
$ws.Range("AW27").Value = "132"
$ws.Range("AX27").NumberFormat = "@"
$ws.Range("AX27").Value = "17.3"
$ws.Range("AY27").NumberFormat = "@"
$ws.Range("AY27").Value = "1022.6"
$ws.Range("AZ27").NumberFormat = "@"
$ws.Range("AZ27").Value = "132"

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "07:30 - 08:00"
$ws.Range("Q32").NumberFormat = "@"
$ws.Range("Q32").Value = "59"
$ws.Range("R32").NumberFormat = "@"
$ws.Range("R32").Value = "1022.7"
$ws.Range("T32").NumberFormat = "@"
$ws.Range("T32").Value = "07:30 - 08:00"
$ws.Range("U32").NumberFormat = "@"
$ws.Range("U32").Value = "120"
$ws.Range("W32").NumberFormat = "@"
$ws.Range("W32").Value = "7.0"
$ws.Range("X32").NumberFormat = "@"
$ws.Range("X32").Value = "7.2"
$ws.Range("Y32").NumberFormat = "@"
$ws.Range("Y32").Value = "49.0"
$ws.Range("AB32").NumberFormat = "@"
$ws.Range("AB32").Value = "57.2"
$ws.Range("AE32").NumberFormat = "@"
$ws.Range("AE32").Value = "07:30 - 08:00"
$ws.Range("AG32").NumberFormat = "@"
$ws.Range("AG32").Value = "7.2"
$ws.Range("AH32").NumberFormat = "@"
$ws.Range("AH32").Value = "7.0"
$ws.Range("AI32").NumberFormat = "@"
$ws.Range("AI32").Value = "59"
$ws.Range("AK32").NumberFormat = "@"
$ws.Range("AK32").Value = "49.0"
$ws.Range("AM32").NumberFormat = "@"
$ws.Range("AM32").Value = "07:30 - 08:00"
$ws.Range("AO32").NumberFormat = "@"
$ws.Range("AO32").Value = "7.2"
$ws.Range("AP32").NumberFormat = "@"
$ws.Range("AP32").Value = "7.0"
$ws.Range("AQ32").NumberFormat = "@"
$ws.Range("AQ32").Value = "59"
$ws.Range("AS32").NumberFormat = "@"
$ws.Range("AS32").Value = "49.0"
$ws.Range("AU32").NumberFormat = "@"
$ws.Range("AU32").Value = "57.2"
$ws.Range("AV32").NumberFormat = "@"
$ws.Range("AV32").Value = "1022.7"
$ws.Range("AW32").NumberFormat = "@"
$ws.Range("AW32").Value = "120"
$ws.Range("AX32").NumberFormat = "@"
$ws.Range("AX32").Value = "57.2"
$ws.Range("AY32").NumberFormat = "@"
$ws.Range("AY32").Value = "1022.7"
$ws.Range("AZ32").NumberFormat = "@"
$ws.Range("AZ32").Value = "120"